$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$scratch = $ws.Range("Z1")
$scratch.NumberFormat = "@"

$ws.Range("D2").Value = "43.048.35"
$ws.Range("E2").Value = "  -0.30%  "

$ws.Range("D3").Value = "2.299.88"
$ws.Range("E3").Value = "  -0.14%  "

$ws.Range("E4").Value = "  +0.03%  "

$scratch.Value = "300.58"
$scratch.Copy()
$ws.Range("D5").PasteSpecial(-4163)
$ws.Range("E5").Value = "  -0.11%  "

$scratch.Value = "99.55"
$scratch.Copy()
$ws.Range("D6").PasteSpecial(-4163)
$ws.Range("E6").Value = "  +2.01%  "

$scratch.Value = "0.505"
$scratch.Copy()
$ws.Range("D7").PasteSpecial(-4163)
$ws.Range("E7").Value = "  -0.82%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("E9").Value = "  +1.81%  "

$scratch.Value = "36.23"
$scratch.Copy()
$ws.Range("D10").PasteSpecial(-4163)
$ws.Range("E10").Value = "  +7.21%  "

$scratch.Value = "0.0790"
$scratch.Copy()
$ws.Range("D11").PasteSpecial(-4163)
$ws.Range("E11").Value = "  -0.85%  "

$ws.Range("E12").Value = "  +0.94%  "

$ws.Range("E13").Value = "  +6.69%  "

$scratch.Value = "6.92"
$scratch.Copy()
$ws.Range("D14").PasteSpecial(-4163)
$ws.Range("E14").Value = "  +1.76%  "

$ws.Range("D15").Value = "2.658.86"
$ws.Range("E15").Value = "  -0.02%  "

$ws.Range("D16").Value = "2.313.47"
$ws.Range("E16").Value = "  +0.48%  "

$ws.Range("E17").Value = "  -1.51%  "

$ws.Range("D18").Value = "42.933.15"
$ws.Range("E18").Value = "  -0.27%  "

$scratch.Value = "12.61"
$scratch.Copy()
$ws.Range("D19").PasteSpecial(-4163)
$ws.Range("E19").Value = "  +8.57%  "

$ws.Range("D20").Value = "0.0₃0905"
$ws.Range("E20").Value = "  +0.21%  "

$scratch.Value = "6.13"
$scratch.Copy()
$ws.Range("D21").PasteSpecial(-4163)
$ws.Range("E21").Value = "  +0.99%  "

$ws.Range("E22").Value = "  +0.26%  "

$scratch.Value = "235.57"
$scratch.Copy()
$ws.Range("D23").PasteSpecial(-4163)
$ws.Range("E23").Value = "  -0.64%  "

$scratch.Value = "2.19"
$scratch.Copy()
$ws.Range("D24").PasteSpecial(-4163)
$ws.Range("E24").Value = "  +8.34%  "

$scratch.Value = "1.00"
$scratch.Copy()
$ws.Range("D25").PasteSpecial(-4163)
$ws.Range("E25").Value = "  +0.50%  "

$ws.Range("E26").Value = "  -0.78%  "

$scratch.Value = "24.94"
$scratch.Copy()
$ws.Range("D27").PasteSpecial(-4163)
$ws.Range("E27").Value = "  +1.45%  "

$ws.Range("E28").Value = "  +15.39%  "

$scratch.Value = "34.49"
$scratch.Copy()
$ws.Range("D29").PasteSpecial(-4163)
$ws.Range("E29").Value = "  +1.08%  "

$scratch.Value = "167.36"
$scratch.Copy()
$ws.Range("D30").PasteSpecial(-4163)
$ws.Range("E30").Value = "  +0.27%  "

$ws.Range("E31").Value = "  -0.22%  "

$ws.Range("E32").Value = "  -0.04%  "

$ws.Range("E33").Value = "  +1.25%  "

$scratch.Value = "17.65"
$scratch.Copy()
$ws.Range("D34").PasteSpecial(-4163)
$ws.Range("E34").Value = "  +3.92%  "

$scratch.Value = "4.58"
$scratch.Copy()
$ws.Range("D35").PasteSpecial(-4163)
$ws.Range("E35").Value = "  -2.41%  "

$scratch.Value = "0.0689"
$scratch.Copy()
$ws.Range("D37").PasteSpecial(-4163)
$ws.Range("E37").Value = "  -0.93%  "

$ws.Range("E38").Value = "  +1.06%  "

$ws.Range("E39").Value = "  -0.78%  "

$ws.Range("E40").Value = "  -0.98%  "

$ws.Range("E41").Value = "  -0.06%  "

$ws.Range("E42").Value = "  +2.98%  "

$ws.Range("E43").Value = "  -3.37%  "

$ws.Range("D44").Value = "1.979.98"
$ws.Range("E44").Value = "  -0.17%  "

$scratch.Value = "10.19"
$scratch.Copy()
$ws.Range("D45").PasteSpecial(-4163)
$ws.Range("E45").Value = "  +2.88%  "

$scratch.Value = "2.91"
$scratch.Copy()
$ws.Range("D46").PasteSpecial(-4163)
$ws.Range("E46").Value = "  +1.37%  "

$scratch.Value = "17.49"
$scratch.Copy()
$ws.Range("D47").PasteSpecial(-4163)
$ws.Range("E47").Value = "  -1.06%  "

$scratch.Value = "55.44"
$scratch.Copy()
$ws.Range("D48").PasteSpecial(-4163)
$ws.Range("E48").Value = "  +3.97%  "

$ws.Range("E49").Value = "  +3.51%  "

$ws.Range("D50").Value = "2.525.20"
$ws.Range("E50").Value = "  -0.05%  "

$scratch.Value = "70.76"
$scratch.Copy()
$ws.Range("D51").PasteSpecial(-4163)
$ws.Range("E51").Value = "  +0.83%  "

$scratch.Clear()
$excel.CutCopyMode = $false
